# The underlying change reverts a schema-serialization tweak (explicit
# default-value attributes like state="visible", outline="1",
# quotePrefix="0", pivotButton="0", customFormat="0", etc. are no longer
# emitted). It does not alter any cell value, sheet name, or formatting,
# so there is nothing to modify through the object model - the workbook
# content and structure stay exactly as they were.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
